$wb = $excel.ActiveWorkbook

# ---- Baja California Mexico ----
$ws = $wb.Worksheets.Item("Baja California Mexico")
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.0161
$ws.Range("O4").Value = 0.00265
$ws.Range("P4").Value = 0.00265
$ws.Range("Q4").Value = 0.00265
$ws.Range("R4").Value = 0.00795
$ws.Range("S4").Value = 0.00265
$ws.Range("T4").Value = 0.00265
$ws.Range("U4").Value = 0.00265
$ws.Range("V4").Value = 0.00795
$ws.Range("W4").Value = 0.0318
$ws.Range("M5").Value = $null
$ws.Range("N5").Value = $null
$ws.Range("M8").Value = 0
$ws.Range("N8").Value = 0
$ws.Range("O8").Value = 0.000583333333333333
$ws.Range("P8").Value = 0.000583333333333333
$ws.Range("Q8").Value = 0.000583333333333333
$ws.Range("R8").Value = 0.00175
$ws.Range("S8").Value = 0.000583333333333333
$ws.Range("T8").Value = 0.000583333333333333
$ws.Range("U8").Value = 0.000583333333333333
$ws.Range("V8").Value = 0.00175
$ws.Range("W8").Value = 0.007

# ---- Cleveland Ohio ----
$ws = $wb.Worksheets.Item("Cleveland Ohio")
$ws.Range("E7").Value = 0.1009
$ws.Range("E8").Value = 0.1009
$ws.Range("E9").Value = 0.1009
$ws.Range("M9").Value = 0
$ws.Range("N9").Value = 0.0327
$ws.Range("O9").Value = 0.0168166666666667
$ws.Range("P9").Value = 0.0168166666666667
$ws.Range("Q9").Value = 0.0168166666666667
$ws.Range("R9").Value = 0.05045
$ws.Range("S9").Value = 0.0168166666666667
$ws.Range("T9").Value = 0.0168166666666667
$ws.Range("U9").Value = 0.0168166666666667
$ws.Range("V9").Value = 0.05045
$ws.Range("W9").Value = 0.2018

# ---- Devon United Kingdom ----
$ws = $wb.Worksheets.Item("Devon United Kingdom")
$ws.Range("E2").Value = 0.0599
$ws.Range("E3").Value = 0.0599
$ws.Range("E4").Value = 0.0599
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.0613
$ws.Range("O4").Value = 0.00998333333333333
$ws.Range("P4").Value = 0.00998333333333333
$ws.Range("Q4").Value = 0.00998333333333333
$ws.Range("R4").Value = 0.02995
$ws.Range("S4").Value = 0.00998333333333333
$ws.Range("T4").Value = 0.00998333333333333
$ws.Range("U4").Value = 0.00998333333333333
$ws.Range("V4").Value = 0.02995
$ws.Range("W4").Value = 0.1198

# ---- Downers Grove Illinois ----
$ws = $wb.Worksheets.Item("Downers Grove Illinois")
$ws.Range("E2").Value = 0.033
$ws.Range("E3").Value = 0.033
$ws.Range("E4").Value = 0.033
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.0055
$ws.Range("P4").Value = 0.0055
$ws.Range("Q4").Value = 0.0055
$ws.Range("R4").Value = 0.0165
$ws.Range("S4").Value = 0.0055
$ws.Range("T4").Value = 0.0055
$ws.Range("U4").Value = 0.0055
$ws.Range("V4").Value = 0.0165
$ws.Range("W4").Value = 0.066
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = 0.6

# ---- East Aurora New York ----
$ws = $wb.Worksheets.Item("East Aurora New York")
$ws.Range("E2").Value = 0.0685
$ws.Range("E3").Value = 0.0685
$ws.Range("E4").Value = 0.0685
$ws.Range("G4").Value = 0.0222
$ws.Range("H4").Value = 0.0227
$ws.Range("J4").Value = 0.0451
$ws.Range("K4").Value = 0.0227
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.0231
$ws.Range("O4").Value = 0.0114166666666667
$ws.Range("P4").Value = 0.0114166666666667
$ws.Range("Q4").Value = 0.0114166666666667
$ws.Range("R4").Value = 0.03425
$ws.Range("S4").Value = 0.0114166666666667
$ws.Range("T4").Value = 0.0114166666666667
$ws.Range("U4").Value = 0.0114166666666667
$ws.Range("V4").Value = 0.03425
$ws.Range("W4").Value = 0.137
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = 0.5
$ws.Range("E8").Value = 0.177
$ws.Range("E9").Value = 0.177
$ws.Range("E10").Value = 0.177
$ws.Range("G10").Value = 0.0417
$ws.Range("H10").Value = 0.0081
$ws.Range("I10").Value = 0.04
$ws.Range("J10").Value = 0.0894
$ws.Range("K10").Value = 0.008
$ws.Range("L10").Value = 0.04
$ws.Range("M10").Value = 0.0394
$ws.Range("N10").Value = 0.0875
$ws.Range("O10").Value = 0.0295
$ws.Range("P10").Value = 0.0295
$ws.Range("Q10").Value = 0.0295
$ws.Range("R10").Value = 0.0885
$ws.Range("S10").Value = 0.0295
$ws.Range("T10").Value = 0.0295
$ws.Range("U10").Value = 0.0295
$ws.Range("V10").Value = 0.0885
$ws.Range("W10").Value = 0.354

# ---- Fremont California ----
$ws = $wb.Worksheets.Item("Fremont California")
$ws.Range("E2").Value = 0.069
$ws.Range("E3").Value = 0.069
$ws.Range("E4").Value = 0.069
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.0115
$ws.Range("P4").Value = 0.0115
$ws.Range("Q4").Value = 0.0115
$ws.Range("R4").Value = 0.0345
$ws.Range("S4").Value = 0.0115
$ws.Range("T4").Value = 0.0115
$ws.Range("U4").Value = 0.0115
$ws.Range("V4").Value = 0.0345
$ws.Range("W4").Value = 0.138
$ws.Range("E5").Value = 0.0383
$ws.Range("E6").Value = 0.0383
$ws.Range("E7").Value = 0.0383
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0.0385
$ws.Range("O7").Value = 0.00638333333333333
$ws.Range("P7").Value = 0.00638333333333333
$ws.Range("Q7").Value = 0.00638333333333333
$ws.Range("R7").Value = 0.01915
$ws.Range("S7").Value = 0.00638333333333333
$ws.Range("T7").Value = 0.00638333333333333
$ws.Range("U7").Value = 0.00638333333333333
$ws.Range("V7").Value = 0.01915
$ws.Range("W7").Value = 0.0766

# ---- Kristianstad Sweden ----
$ws = $wb.Worksheets.Item("Kristianstad Sweden")
$ws.Range("E2").Value = 0.0485
$ws.Range("E3").Value = 0.0485
$ws.Range("E4").Value = 0.0485
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("O4").Value = 0.00808333333333333
$ws.Range("P4").Value = 0.00808333333333333
$ws.Range("Q4").Value = 0.00808333333333333
$ws.Range("R4").Value = 0.02425
$ws.Range("S4").Value = 0.00808333333333333
$ws.Range("T4").Value = 0.00808333333333333
$ws.Range("U4").Value = 0.00808333333333333
$ws.Range("V4").Value = 0.02425
$ws.Range("W4").Value = 0.097
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = $null
$ws.Range("E8").Value = 0.0521
$ws.Range("E9").Value = 0.0521
$ws.Range("E10").Value = 0.0521
$ws.Range("H10").Value = 0.0104
$ws.Range("I10").Value = 0.0103
$ws.Range("J10").Value = 0.0208
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0.0104
$ws.Range("M10").Value = 0.0213
$ws.Range("N10").Value = 0.0313
$ws.Range("O10").Value = 0.00868333333333333
$ws.Range("P10").Value = 0.00868333333333333
$ws.Range("Q10").Value = 0.00868333333333333
$ws.Range("R10").Value = 0.02605
$ws.Range("S10").Value = 0.00868333333333333
$ws.Range("T10").Value = 0.00868333333333333
$ws.Range("U10").Value = 0.00868333333333333
$ws.Range("V10").Value = 0.02605
$ws.Range("W10").Value = 0.1042

# ---- Marengo Illinois ----
$ws = $wb.Worksheets.Item("Marengo Illinois")
$ws.Range("E2").Value = 0.0896
$ws.Range("E3").Value = 0.0896
$ws.Range("E4").Value = 0.0896
$ws.Range("G4").Value = 0.0536
$ws.Range("J4").Value = 0.0536
$ws.Range("K4").Value = 0.0364
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.0359
$ws.Range("O4").Value = 0.0149333333333333
$ws.Range("P4").Value = 0.0149333333333333
$ws.Range("Q4").Value = 0.0149333333333333
$ws.Range("R4").Value = 0.0448
$ws.Range("S4").Value = 0.0149333333333333
$ws.Range("T4").Value = 0.0149333333333333
$ws.Range("U4").Value = 0.0149333333333333
$ws.Range("V4").Value = 0.0448
$ws.Range("W4").Value = 0.1792
$ws.Range("E7").Value = 0.0792
$ws.Range("E8").Value = 0.0792
$ws.Range("E9").Value = 0.0792
$ws.Range("L9").Value = 0.0203
$ws.Range("M9").Value = 0.0138
$ws.Range("N9").Value = 0.0474
$ws.Range("O9").Value = 0.0132
$ws.Range("P9").Value = 0.0132
$ws.Range("Q9").Value = 0.0132
$ws.Range("R9").Value = 0.0396
$ws.Range("S9").Value = 0.0132
$ws.Range("T9").Value = 0.0132
$ws.Range("U9").Value = 0.0132
$ws.Range("V9").Value = 0.0396
$ws.Range("W9").Value = 0.1584

# ---- Milwaukee Pmc Hq Wisconsin ----
$ws = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")
$ws.Range("M5").Value = $null
$ws.Range("N5").Value = $null

# ---- Rock Road Radford Virginia ----
$ws = $wb.Worksheets.Item("Rock Road Radford Virginia")
$ws.Range("E2").Value = 0.0157
$ws.Range("E3").Value = 0.0157
$ws.Range("E4").Value = 0.0157
$ws.Range("L4").Value = 0.0159
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0.0158
$ws.Range("O4").Value = 0.00261666666666667
$ws.Range("P4").Value = 0.00261666666666667
$ws.Range("Q4").Value = 0.00261666666666667
$ws.Range("R4").Value = 0.00785
$ws.Range("S4").Value = 0.00261666666666667
$ws.Range("T4").Value = 0.00261666666666667
$ws.Range("U4").Value = 0.00261666666666667
$ws.Range("V4").Value = 0.00785
$ws.Range("W4").Value = 0.0314
$ws.Range("M7").Value = $null
$ws.Range("N7").Value = $null

# ---- Shanghai Minhang District Chin ----
$ws = $wb.Worksheets.Item("Shanghai Minhang District Chin")
$ws.Range("E2").Value = 0.2105
$ws.Range("E3").Value = 0.2105
$ws.Range("E4").Value = 0.2105
$ws.Range("M4").Value = 0.1111
$ws.Range("N4").Value = 0.1075
$ws.Range("O4").Value = 0.0350833333333333
$ws.Range("P4").Value = 0.0350833333333333
$ws.Range("Q4").Value = 0.0350833333333333
$ws.Range("R4").Value = 0.10525
$ws.Range("S4").Value = 0.0350833333333333
$ws.Range("T4").Value = 0.0350833333333333
$ws.Range("U4").Value = 0.0350833333333333
$ws.Range("V4").Value = 0.10525
$ws.Range("W4").Value = 0.421

# ---- Betzdorf Germany ----
$ws = $wb.Worksheets.Item("Betzdorf Germany")
$ws.Range("M3").Value = $null
